$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7550.8335
$ws.Range("I19").Value = 2564.3076
$ws.Range("J19").Value = 13444
$ws.Range("K19").Value = 2564.3076
$ws.Range("L19").Value = 13444
$ws.Range("M19").Value = -2389.3076
$ws.Range("N19").Value = -13794

$ws.Range("H20").Value = 2999
$ws.Range("I20").Value = 2999
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2999
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2769
$ws.Range("N20").ClearContents()

$ws.Range("H21").Value = 750
$ws.Range("I21").Value = 750
$ws.Range("K21").Value = 750
$ws.Range("M21").Value = -282

$ws.Range("H23").Value = 750
$ws.Range("I23").Value = 750
$ws.Range("K23").Value = 750
$ws.Range("M23").Value = -516

$ws.Range("H31").Value = 23666.666
$ws.Range("I31").Value = 23666.666
$ws.Range("K31").Value = 70999.99800000001
$ws.Range("M31").Value = -70769.99800000001

$ws.Range("H35").Value = 2999
$ws.Range("I35").Value = 2999
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2999
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2620
$ws.Range("N35").ClearContents()

$ws.Range("H38").Value = 3184.25
$ws.Range("I38").Value = 113.2
$ws.Range("J38").Value = 4207.933
$ws.Range("K38").Value = 339.6
$ws.Range("L38").Value = 12623.799
$ws.Range("M38").Value = 32.39999999999998
$ws.Range("N38").Value = -13367.799

$ws.Range("H96").Value = 52634580
$ws.Range("I96").Value = 3899.1538
$ws.Range("J96").Value = 166667710
$ws.Range("K96").Value = 11697.4614
$ws.Range("L96").Value = 500003130
$ws.Range("M96").Value = -10324.4614
$ws.Range("N96").Value = -500005876

$ws.Range("H112").Value = 2453.4167
$ws.Range("J112").Value = 3241.1667
$ws.Range("L112").Value = 9723.500100000001
$ws.Range("N112").Value = -11939.5001

$ws.Range("H125").Value = 3774.5652
$ws.Range("I125").Value = 1647.4286
$ws.Range("J125").Value = 7083.4443
$ws.Range("K125").Value = 14826.8574
$ws.Range("L125").Value = 63750.9987
$ws.Range("M125").Value = -12366.8574
$ws.Range("N125").Value = -68670.9987

$ws.Range("H132").Value = 6148.491
$ws.Range("I132").Value = 3822.7144
$ws.Range("J132").Value = 12660.667
$ws.Range("K132").Value = 11468.1432
$ws.Range("L132").Value = 37982.001
$ws.Range("M132").Value = -8938.143199999999
$ws.Range("N132").Value = -43042.001

$ws.Range("H138").Value = 2326.09
$ws.Range("I138").Value = 1600.6471
$ws.Range("J138").Value = 2699.803
$ws.Range("K138").Value = 4801.9413
$ws.Range("L138").Value = 8099.409
$ws.Range("M138").Value = 338.0587000000005
$ws.Range("N138").Value = -18379.409

$ws.Range("H141").Value = 10706.538
$ws.Range("I141").Value = 12491.346
$ws.Range("K141").Value = 37474.038
$ws.Range("M141").Value = -32294.038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 28166.166
$ws.Range("J96").Value = 28166.166
$ws.Range("L96").Value = 28166.166
$ws.Range("N96").Value = -33658.166

$ws.Range("H122").Value = 4432
$ws.Range("I122").Value = 2222
$ws.Range("J122").Value = 4984.5
$ws.Range("K122").Value = 6666
$ws.Range("L122").Value = 14953.5
$ws.Range("M122").Value = -4216
$ws.Range("N122").Value = -19853.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3882.1428
$ws.Range("I86").Value = 2793
$ws.Range("J86").Value = 4699
$ws.Range("K86").Value = 2793
$ws.Range("L86").Value = 4699
$ws.Range("M86").Value = -1670
$ws.Range("N86").Value = -6945

$ws.Range("H89").Value = 3882.1428
$ws.Range("I89").Value = 2793
$ws.Range("J89").Value = 4699
$ws.Range("K89").Value = 13965
$ws.Range("L89").Value = 23495
$ws.Range("M89").Value = -8349
$ws.Range("N89").Value = -34727

$ws.Range("H94").Value = 1261.5278
$ws.Range("I94").Value = 762.8077
$ws.Range("K94").Value = 762.8077
$ws.Range("M94").Value = -311.8077

$ws.Range("H99").Value = 1966.35
$ws.Range("I99").Value = 1852.0555
$ws.Range("K99").Value = 1852.0555
$ws.Range("M99").Value = -354.0554999999999

$ws.Range("H134").Value = 2096.9028
$ws.Range("I134").Value = 1405.2034
$ws.Range("J134").Value = 5236.154
$ws.Range("K134").Value = 4215.6102
$ws.Range("L134").Value = 15708.462
$ws.Range("M134").Value = -1680.6102
$ws.Range("N134").Value = -20778.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12603.963
$ws.Range("I99").Value = 12123.125
$ws.Range("J99").Value = 12806.421
$ws.Range("K99").Value = 12123.125
$ws.Range("L99").Value = 12806.421
$ws.Range("M99").Value = -10625.125
$ws.Range("N99").Value = -15802.421

$ws.Range("H111").Value = 130633
$ws.Range("J111").Value = 130633
$ws.Range("L111").Value = 130633
$ws.Range("N111").Value = -138813

$ws.Range("H122").Value = 101800
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 101800
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 305400
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -310300

$ws.Range("H126").Value = 12603.963
$ws.Range("I126").Value = 12123.125
$ws.Range("J126").Value = 12806.421
$ws.Range("K126").Value = 36369.375
$ws.Range("L126").Value = 38419.263
$ws.Range("M126").Value = -33899.375
$ws.Range("N126").Value = -43359.263

$ws.Range("H132").Value = 2817.718
$ws.Range("I132").Value = 2174.2258
$ws.Range("J132").Value = 5311.25
$ws.Range("K132").Value = 6522.6774
$ws.Range("L132").Value = 15933.75
$ws.Range("M132").Value = -3992.6774
$ws.Range("N132").Value = -20993.75

$ws.Range("H134").Value = 3330.2036
$ws.Range("I134").Value = 2833.6047
$ws.Range("K134").Value = 8500.8141
$ws.Range("M134").Value = -5965.8141

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H113").Value = 1316.5385
$ws.Range("I113").Value = 1441.6
$ws.Range("J113").Value = 899.6667
$ws.Range("K113").Value = 4324.799999999999
$ws.Range("L113").Value = 2699.0001
$ws.Range("M113").Value = -2154.799999999999
$ws.Range("N113").Value = -7039.0001

$ws.Range("H137").Value = 5010
$ws.Range("J137").Value = 8411
$ws.Range("L137").Value = 25233
$ws.Range("N137").Value = -35433

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6483
$ws.Range("I102").Value = 5678
$ws.Range("K102").Value = 5678
$ws.Range("M102").Value = -4056

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 71431140
$ws.Range("I100").Value = 166668640
$ws.Range("K100").Value = 166668640
$ws.Range("M100").Value = -166668099

$ws.Range("H122").Value = 6872.643
$ws.Range("I122").Value = 7087.391
$ws.Range("K122").Value = 21262.173
$ws.Range("M122").Value = -18812.173

$ws.Range("H132").Value = 31919.59
$ws.Range("I132").Value = 51794.273
$ws.Range("J132").Value = 6199.4116
$ws.Range("K132").Value = 155382.819
$ws.Range("L132").Value = 18598.2348
$ws.Range("M132").Value = -152852.819
$ws.Range("N132").Value = -23658.2348

$ws.Range("H136").Value = 6930352
$ws.Range("I136").Value = 8189384.5
$ws.Range("K136").Value = 24568153.5
$ws.Range("M136").Value = -24565603.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 13245
$ws.Range("J8").Value = 15993.333
$ws.Range("L8").Value = 15993.333
$ws.Range("N8").Value = -16273.333

$ws.Range("H122").Value = 6034.022
$ws.Range("I122").Value = 3727.8333
$ws.Range("K122").Value = 11183.4999
$ws.Range("M122").Value = -8733.499899999999

Write-Host "Applied all changes"